# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H, I, J, K, L, M, N) across several leve rows on
# several sheets with newly-fetched values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 5005.6665
$ws.Range("I15").Value = 5005.6665
$ws.Range("K15").Value = 15016.9995
$ws.Range("M15").Value = -14847.9995

$ws.Range("H28").Value = 590.76
$ws.Range("I28").Value = 552.9167
$ws.Range("K28").Value = 552.9167
$ws.Range("M28").Value = -67.91669999999999

$ws.Range("H64").Value = 4050
$ws.Range("I64").Value = 3600
$ws.Range("K64").Value = 3600
$ws.Range("M64").Value = -3352

$ws.Range("H67").Value = 4050
$ws.Range("I67").Value = 3600
$ws.Range("K67").Value = 3600
$ws.Range("M67").Value = -2742

$ws.Range("H103").Value = 1419.125
$ws.Range("I103").Value = 1016.2857
$ws.Range("K103").Value = 3048.8571
$ws.Range("M103").Value = -2462.8571

$ws.Range("H132").Value = 4556.029
$ws.Range("I132").Value = 2400.4333
$ws.Range("K132").Value = 7201.2999
$ws.Range("M132").Value = -4671.2999

$ws.Range("H137").Value = 8607.267
$ws.Range("J137").Value = 4994.5
$ws.Range("L137").Value = 14983.5
$ws.Range("N137").Value = -20083.5

$ws.Range("H138").Value = 2655.0476
$ws.Range("J138").Value = 2683.7144
$ws.Range("L138").Value = 8051.1432
$ws.Range("N138").Value = -18331.1432

$ws.Range("H140").Value = 69392
$ws.Range("J140").Value = 69392
$ws.Range("L140").Value = 69392
$ws.Range("N140").Value = -79752

$ws.Range("H141").Value = 11733.92
$ws.Range("I141").Value = 11232.956
$ws.Range("K141").Value = 33698.868
$ws.Range("M141").Value = -28518.868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6406.732
$ws.Range("I32").Value = 6032.9326
$ws.Range("K32").Value = 6032.9326
$ws.Range("M32").Value = -5745.9326

$ws.Range("H74").Value = 1956.6389
$ws.Range("I74").Value = 1450.9
$ws.Range("K74").Value = 1450.9
$ws.Range("M74").Value = -576.9000000000001

$ws.Range("H77").Value = 1956.6389
$ws.Range("I77").Value = 1450.9
$ws.Range("K77").Value = 7254.5
$ws.Range("M77").Value = -2886.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 26651.2
$ws.Range("I82").Value = 7085.6665
$ws.Range("J82").Value = 55999.5
$ws.Range("K82").Value = 7085.6665
$ws.Range("L82").Value = 55999.5
$ws.Range("M82").Value = -6702.6665
$ws.Range("N82").Value = -56765.5

$ws.Range("H85").Value = 26651.2
$ws.Range("I85").Value = 7085.6665
$ws.Range("J85").Value = 55999.5
$ws.Range("K85").Value = 7085.6665
$ws.Range("L85").Value = 55999.5
$ws.Range("M85").Value = -5759.6665
$ws.Range("N85").Value = -58651.5

$ws.Range("H86").Value = 3391.5
$ws.Range("I86").Value = 3337.5
$ws.Range("J86").Value = 3499.5
$ws.Range("K86").Value = 3337.5
$ws.Range("L86").Value = 3499.5
$ws.Range("M86").Value = -2214.5
$ws.Range("N86").Value = -5745.5

$ws.Range("H89").Value = 3391.5
$ws.Range("I89").Value = 3337.5
$ws.Range("J89").Value = 3499.5
$ws.Range("K89").Value = 16687.5
$ws.Range("L89").Value = 17497.5
$ws.Range("M89").Value = -11071.5
$ws.Range("N89").Value = -28729.5

$ws.Range("H99").Value = 3610.1538
$ws.Range("I99").Value = 1784.2667
$ws.Range("K99").Value = 1784.2667
$ws.Range("M99").Value = -286.2666999999999

$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3149.4722
$ws.Range("I31").Value = 1922.2
$ws.Range("K31").Value = 1922.2
$ws.Range("M31").Value = -1627.2

$ws.Range("H34").Value = 3149.4722
$ws.Range("I34").Value = 1922.2
$ws.Range("K34").Value = 1922.2
$ws.Range("M34").Value = -1720.2

$ws.Range("H58").Value = 4845.4375
$ws.Range("I58").Value = 4445.154
$ws.Range("K58").Value = 4445.154
$ws.Range("M58").Value = -4242.154

$ws.Range("H132").Value = 10461.63
$ws.Range("I132").Value = 3960.7856
$ws.Range("K132").Value = 11882.3568
$ws.Range("M132").Value = -9352.356800000001

$ws.Range("H136").Value = 4845.4375
$ws.Range("I136").Value = 4445.154
$ws.Range("K136").Value = 13335.462
$ws.Range("M136").Value = -10785.462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 739
$ws.Range("J92").Value = 724.8570999999999
$ws.Range("L92").Value = 2174.5713
$ws.Range("N92").Value = -4670.5713

$ws.Range("H97").Value = 5109.857
$ws.Range("I97").Value = 4744
$ws.Range("J97").Value = 5597.6665
$ws.Range("K97").Value = 14232
$ws.Range("L97").Value = 16792.9995
$ws.Range("M97").Value = -13736
$ws.Range("N97").Value = -17784.9995

$ws.Range("H98").Value = 2501561.5
$ws.Range("I98").Value = 5001749.5
$ws.Range("K98").Value = 15005248.5
$ws.Range("M98").Value = -15003750.5

$ws.Range("H109").Value = 1513.2858
$ws.Range("I109").Value = 1429.6666
$ws.Range("K109").Value = 4288.9998
$ws.Range("M109").Value = -3248.9998

$ws.Range("H122").Value = 1375.5834
$ws.Range("J122").Value = 1555.8
$ws.Range("L122").Value = 14002.2
$ws.Range("N122").Value = -18902.2

$ws.Range("H139").Value = 2798.8096
$ws.Range("I139").Value = 1777.6316
$ws.Range("J139").Value = 12500
$ws.Range("K139").Value = 5332.8948
$ws.Range("L139").Value = 37500
$ws.Range("M139").Value = -192.8948
$ws.Range("N139").Value = -47780

$ws.Range("H140").Value = 949.2632
$ws.Range("I140").Value = 949.2632
$ws.Range("K140").Value = 2847.7896
$ws.Range("M140").Value = 2332.2104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 32866300
$ws.Range("J80").Value = 2998.5
$ws.Range("L80").Value = 2998.5
$ws.Range("N80").Value = -4994.5

$ws.Range("H83").Value = 32866300
$ws.Range("J83").Value = 2998.5
$ws.Range("L83").Value = 14992.5
$ws.Range("N83").Value = -24976.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2758.4546
$ws.Range("I7").Value = 2234.3
$ws.Range("K7").Value = 2234.3
$ws.Range("M7").Value = -2122.3

$ws.Range("H68").Value = 2641.4
$ws.Range("I68").Value = 2693.7693
$ws.Range("J68").Value = 2301
$ws.Range("K68").Value = 2693.7693
$ws.Range("L68").Value = 2301
$ws.Range("M68").Value = -1944.7693
$ws.Range("N68").Value = -3799

$ws.Range("H71").Value = 2641.4
$ws.Range("I71").Value = 2693.7693
$ws.Range("J71").Value = 2301
$ws.Range("K71").Value = 13468.8465
$ws.Range("L71").Value = 11505
$ws.Range("M71").Value = -9724.8465
$ws.Range("N71").Value = -18993

$ws.Range("H110").Value = 64469
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 64469
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 64469
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -72649

$ws.Range("H126").Value = 2758.4546
$ws.Range("I126").Value = 2234.3
$ws.Range("K126").Value = 6702.900000000001
$ws.Range("M126").Value = -4232.900000000001

$ws.Range("H132").Value = 18122.027
$ws.Range("I132").Value = 21070.89
$ws.Range("K132").Value = 63212.67
$ws.Range("M132").Value = -60682.67

$ws.Range("H139").Value = 120268.86
$ws.Range("J139").Value = 108776.4
$ws.Range("L139").Value = 108776.4
$ws.Range("N139").Value = -119056.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 150000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 150000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 150000
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -152774

$ws.Range("H132").Value = 4269.365
$ws.Range("I132").Value = 4083.72
$ws.Range("K132").Value = 12251.16
$ws.Range("M132").Value = -9721.16
